$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data (columns D, J, K, L, M, O, P) for rows 2..12,
# representing a re-ordering of the existing records by date.
$data = @(
    @{ Row = 2;  D = 44637; J = 170; K = 2800; L = 3000; M = 2906; O = "Región Metropolitana";   P = 484 }
    @{ Row = 3;  D = 44643; J = 90;  K = 2800; L = 3000; M = 2911; O = "Región Metropolitana";   P = 485 }
    @{ Row = 4;  D = 44658; J = 180; K = 2500; L = 3000; M = 2778; O = "Región Metropolitana";   P = 463 }
    @{ Row = 5;  D = 44659; J = 90;  K = 2500; L = 3000; M = 2722; O = "Región Metropolitana";   P = 454 }
    @{ Row = 6;  D = 44631; J = 110; K = 3000; L = 3500; M = 3273; O = "Provincia de Chacabuco"; P = 546 }
    @{ Row = 7;  D = 44644; J = 140; K = 2500; L = 3000; M = 2786; O = "Provincia de Chacabuco"; P = 464 }
    @{ Row = 8;  D = 44672; J = 140; K = 3000; L = 3500; M = 3286; O = "Región Metropolitana";   P = 548 }
    @{ Row = 9;  D = 44671; J = 150; K = 3500; L = 4000; M = 3733; O = "Región Metropolitana";   P = 622 }
    @{ Row = 10; D = 44685; J = 150; K = 3000; L = 3500; M = 3267; O = "Región Metropolitana";   P = 544 }
    @{ Row = 11; D = 44630; J = 90;  K = 2500; L = 3000; M = 2722; O = "Región Metropolitana";   P = 454 }
    @{ Row = 12; D = 44650; J = 130; K = 3000; L = 3500; M = 3308; O = "Región Metropolitana";   P = 551 }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 4).Value  = $rec.D   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $rec.J   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $rec.K   # K - Precio mínimo
    $ws.Cells.Item($r, 12).Value = $rec.L   # L - Precio máximo
    $ws.Cells.Item($r, 13).Value = $rec.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value = $rec.O   # O - Origen
    $ws.Cells.Item($r, 16).Value = $rec.P   # P - Precio $/Kg
}
